# Auto-generated edit script applying value updates per the commit diff.
# Updates numeric cells (H:N) across multiple sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 2021.5
$ws.Range("I62").Value = 1912.7778
$ws.Range("K62").Value = 1912.7778
$ws.Range("M62").Value = -1288.7778

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 2021.5
$ws.Range("I65").Value = 1912.7778
$ws.Range("K65").Value = 9563.889000000001
$ws.Range("M65").Value = -6443.889000000001

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 61377.117
$ws.Range("I86").Value = 201540.6
$ws.Range("J86").Value = 2975.6667
$ws.Range("K86").Value = 201540.6
$ws.Range("L86").Value = 2975.6667
$ws.Range("M86").Value = -200417.6
$ws.Range("N86").Value = -5221.6667

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 61377.117
$ws.Range("I89").Value = 201540.6
$ws.Range("J89").Value = 2975.6667
$ws.Range("K89").Value = 1007703
$ws.Range("L89").Value = 14878.3335
$ws.Range("M89").Value = -1002087
$ws.Range("N89").Value = -26110.3335

# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 5191.6816
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 5191.6816
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 46725.1344
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -51645.1344

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 2715.1667
$ws.Range("I2").Value = 2715.1667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2715.1667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2602.1667
$ws.Range("N2").ClearContents()

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1370.2333
$ws.Range("I45").Value = 1299.8148
$ws.Range("J45").Value = 2004
$ws.Range("K45").Value = 1299.8148
$ws.Range("L45").Value = 2004
$ws.Range("M45").Value = -922.8148000000001
$ws.Range("N45").Value = -2758

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 6851.1133
$ws.Range("I61").Value = 5734.575
$ws.Range("J61").Value = 10286.615
$ws.Range("K61").Value = 5734.575
$ws.Range("L61").Value = 10286.615
$ws.Range("M61").Value = -5522.575
$ws.Range("N61").Value = -10710.615

# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 1341.1
$ws.Range("I97").Value = 1314.2858
$ws.Range("J97").Value = 1403.6666
$ws.Range("K97").Value = 1314.2858
$ws.Range("L97").Value = 1403.6666
$ws.Range("M97").Value = -818.2858000000001
$ws.Range("N97").Value = -2395.6666

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 1405.8572
$ws.Range("I110").Value = 1358.3334
$ws.Range("J110").Value = 1691
$ws.Range("K110").Value = 1358.3334
$ws.Range("L110").Value = 1691
$ws.Range("M110").Value = 686.6666
$ws.Range("N110").Value = -5781

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 2715.1667
$ws.Range("I116").Value = 2715.1667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2715.1667
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -421.1667000000002
$ws.Range("N116").ClearContents()

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3963.6316
$ws.Range("I132").Value = 1565.3243
$ws.Range("J132").Value = 8400.5
$ws.Range("K132").Value = 4695.9729
$ws.Range("L132").Value = 25201.5
$ws.Range("M132").Value = -2165.9729
$ws.Range("N132").Value = -30261.5

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 6851.1133
$ws.Range("I136").Value = 5734.575
$ws.Range("J136").Value = 10286.615
$ws.Range("K136").Value = 17203.725
$ws.Range("L136").Value = 30859.845
$ws.Range("M136").Value = -14653.725
$ws.Range("N136").Value = -35959.845

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 2715.1667
$ws.Range("I3").Value = 2715.1667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2715.1667
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2601.1667
$ws.Range("N3").ClearContents()

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 1916.2162
$ws.Range("I86").Value = 1950
$ws.Range("J86").Value = 1533.3334
$ws.Range("K86").Value = 1950
$ws.Range("L86").Value = 1533.3334
$ws.Range("M86").Value = -827
$ws.Range("N86").Value = -3779.3334

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 1916.2162
$ws.Range("I89").Value = 1950
$ws.Range("J89").Value = 1533.3334
$ws.Range("K89").Value = 9750
$ws.Range("L89").Value = 7666.666999999999
$ws.Range("M89").Value = -4134
$ws.Range("N89").Value = -18898.667

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 995.2353000000001
$ws.Range("I94").Value = 989.0833
$ws.Range("J94").Value = 1010
$ws.Range("K94").Value = 989.0833
$ws.Range("L94").Value = 1010
$ws.Range("M94").Value = -538.0833
$ws.Range("N94").Value = -1912

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 804286.8
$ws.Range("I105").Value = 1159736.5
$ws.Range("J105").Value = 4525
$ws.Range("K105").Value = 1159736.5
$ws.Range("L105").Value = 4525
$ws.Range("M105").Value = -1157989.5
$ws.Range("N105").Value = -8019

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 1181.2916
$ws.Range("I107").Value = 920.8
$ws.Range("J107").Value = 1367.3572
$ws.Range("K107").Value = 920.8
$ws.Range("L107").Value = 1367.3572
$ws.Range("M107").Value = 999.2
$ws.Range("N107").Value = -5207.3572

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 6751.0356
$ws.Range("I31").Value = 9861.154
$ws.Range("K31").Value = 9861.154
$ws.Range("M31").Value = -9566.154

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 6751.0356
$ws.Range("I34").Value = 9861.154
$ws.Range("K34").Value = 9861.154
$ws.Range("M34").Value = -9659.154

# Row 60 (Leve Item ID 1937)
$ws.Range("H60").Value = 29000
$ws.Range("J60").Value = 29000
$ws.Range("L60").Value = 29000
$ws.Range("N60").Value = -30022

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 3062.5
$ws.Range("I99").Value = 2320
$ws.Range("J99").Value = 4300
$ws.Range("K99").Value = 2320
$ws.Range("L99").Value = 4300
$ws.Range("M99").Value = -822
$ws.Range("N99").Value = -7296

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 3062.5
$ws.Range("I126").Value = 2320
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 6960
$ws.Range("L126").Value = 12900
$ws.Range("M126").Value = -4490
$ws.Range("N126").Value = -17840

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 10267.286
$ws.Range("I122").Value = 12589.556
$ws.Range("J122").Value = 8525.583000000001
$ws.Range("K122").Value = 37768.66800000001
$ws.Range("L122").Value = 25576.749
$ws.Range("M122").Value = -35318.66800000001
$ws.Range("N122").Value = -30476.749

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2355.9285
$ws.Range("I126").Value = 1587.3334
$ws.Range("J126").Value = 2932.375
$ws.Range("K126").Value = 4762.0002
$ws.Range("L126").Value = 8797.125
$ws.Range("M126").Value = -2292.0002
$ws.Range("N126").Value = -13737.125

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 851.53845
$ws.Range("I22").Value = 833.5
$ws.Range("J22").Value = 862.8125
$ws.Range("K22").Value = 833.5
$ws.Range("L22").Value = 862.8125
$ws.Range("M22").Value = -538.5
$ws.Range("N22").Value = -1452.8125

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 851.53845
$ws.Range("I27").Value = 833.5
$ws.Range("J27").Value = 862.8125
$ws.Range("K27").Value = 833.5
$ws.Range("L27").Value = 862.8125
$ws.Range("M27").Value = -726.5
$ws.Range("N27").Value = -1076.8125

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1233.6666
$ws.Range("I46").Value = 2250.5
$ws.Range("J46").Value = 943.1429000000001
$ws.Range("K46").Value = 2250.5
$ws.Range("L46").Value = 943.1429000000001
$ws.Range("M46").Value = -2062.5
$ws.Range("N46").Value = -1319.1429

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 299
$ws.Range("J55").Value = 372
$ws.Range("L55").Value = 372
$ws.Range("N55").Value = -718

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 6712.2954
$ws.Range("I132").Value = 8484.593000000001
$ws.Range("J132").Value = 3897.4707
$ws.Range("K132").Value = 25453.779
$ws.Range("L132").Value = 11692.4121
$ws.Range("M132").Value = -22923.779
$ws.Range("N132").Value = -16752.4121

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 1329.9231
$ws.Range("I107").Value = 445
$ws.Range("J107").Value = 1723.2222
$ws.Range("K107").Value = 1335
$ws.Range("L107").Value = 5169.6666
$ws.Range("M107").Value = 585
$ws.Range("N107").Value = -9009.6666

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 1405.6364
$ws.Range("I126").Value = 1570.3846
$ws.Range("J126").Value = 1167.6666
$ws.Range("K126").Value = 4711.1538
$ws.Range("L126").Value = 3502.9998
$ws.Range("M126").Value = -2241.1538
$ws.Range("N126").Value = -8442.9998

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1347.2463
$ws.Range("I132").Value = 652.49054
$ws.Range("J132").Value = 3648.625
$ws.Range("K132").Value = 1957.47162
$ws.Range("L132").Value = 10945.875
$ws.Range("M132").Value = 572.52838
$ws.Range("N132").Value = -16005.875
